$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.331.10'
$ws.Range("D3").Value = '1.931.41'
$ws.Range("E3").Value = '  -2.54%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = '''240.98'
$ws.Range("E5").Value = '  -1.30%  '
$ws.Range("D6").Value = '''0.606'
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("D8").Value = '''56.77'
$ws.Range("E8").Value = '  -4.28%  '
$ws.Range("E9").Value = '  -4.65%  '
$ws.Range("D10").Value = '''0.0833'
$ws.Range("E10").Value = '  +0.91%  '
$ws.Range("E11").Value = '  -0.23%  '
$ws.Range("D12").Value = '2.216.21'
$ws.Range("E12").Value = '  -2.52%  '
$ws.Range("E13").Value = '  -7.14%  '
$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D14").Value = '''13.37'
$ws.Range("E14").Value = '  -4.28%  '
$ws.Range("B15").Value = 'Avalanche'
$ws.Range("C15").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D15").Value = '''20.90'
$ws.Range("E15").Value = '  -11.30%  '
$ws.Range("D16").Value = '''5.13'
$ws.Range("E16").Value = '  -6.00%  '
$ws.Range("D17").Value = '1.953.43'
$ws.Range("E17").Value = '  -1.87%  '
$ws.Range("D18").Value = '36.260.30'
$ws.Range("E18").Value = '  -0.03%  '
$ws.Range("D19").Value = '''68.80'
$ws.Range("E19").Value = '  -1.45%  '
$ws.Range("D20").Value = '0.0₃0858'
$ws.Range("E20").Value = '  -1.00%  '
$ws.Range("D21").Value = '''225.65'
$ws.Range("E21").Value = '  -3.62%  '
$ws.Range("D22").Value = '''4.94'
$ws.Range("E22").Value = '  -6.92%  '
$ws.Range("E23").Value = '  -0.10%  '
$ws.Range("D24").Value = '''2.33'
$ws.Range("E24").Value = '  -10.56%  '
$ws.Range("D25").Value = '''2.26'
$ws.Range("E25").Value = '  -1.93%  '
$ws.Range("D26").Value = '''9.27'
$ws.Range("E26").Value = '  -7.93%  '
$ws.Range("D27").Value = '''160.33'
$ws.Range("E27").Value = '  -1.06%  '
$ws.Range("D28").Value = '''0.130'
$ws.Range("E28").Value = '  -1.19%  '
$ws.Range("D29").Value = '''19.09'
$ws.Range("E29").Value = '  -3.72%  '
$ws.Range("E30").Value = '  -2.49%  '
$ws.Range("E31").Value = '  -6.96%  '
$ws.Range("D32").Value = '''4.53'
$ws.Range("E32").Value = '  -7.44%  '
$ws.Range("D33").Value = '''0.0625'
$ws.Range("E33").Value = '  -0.22%  '
$ws.Range("E34").Value = '  -6.59%  '
$ws.Range("E35").Value = '  -0.01%  '
$ws.Range("D36").Value = '''6.09'
$ws.Range("E36").Value = '  -2.02%  '
$ws.Range("E37").Value = '  -0.15%  '
$ws.Range("D38").Value = '''2.11'
$ws.Range("E38").Value = '  -6.78%  '
$ws.Range("D39").Value = '''2.96'
$ws.Range("E39").Value = '  -2.49%  '
$ws.Range("E40").Value = '  +0.44%  '
$ws.Range("E41").Value = '  -1.34%  '
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").Value = '''1.14'
$ws.Range("E42").Value = '  -7.73%  '
$ws.Range("B43").Value = 'VeChain'
$ws.Range("C43").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D43").Value = '''0.0207'
$ws.Range("E43").Value = '  -3.05%  '
$ws.Range("D44").Value = '''15.45'
$ws.Range("E44").Value = '  -4.85%  '
$ws.Range("D45").Value = '1.331.92'
$ws.Range("E45").Value = '  -2.95%  '
$ws.Range("E46").Value = '  -7.58%  '
$ws.Range("D47").Value = '''86.07'
$ws.Range("E47").Value = '  -6.91%  '
$ws.Range("D48").Value = '''7.05'
$ws.Range("E48").Value = '  -6.02%  '
$ws.Range("E49").Value = '  -0.65%  '
$ws.Range("D50").Value = '''43.78'
$ws.Range("E50").Value = '  -3.52%  '
$ws.Range("D51").Value = '2.108.18'
$ws.Range("E51").Value = '  -2.58%  '
